$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the new (blank) columns at the positions required by the target
#    layout. Each Insert() call shifts everything at/after that column one
#    slot to the right, carrying formulas/values/styles with it.
#    (Column letters below are the CURRENT sheet positions at the time each
#    statement runs, i.e. they already account for the previous inserts.)
# ---------------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Insert()        # new: C
$ws.Range("G1:I1").EntireColumn.Insert()     # new: G, H, I
$ws.Range("M1:P1").EntireColumn.Insert()     # new: M, N, O, P
$ws.Range("R1:V1").EntireColumn.Insert()     # new: R, S, T, U, V

# ---------------------------------------------------------------------------
# 2. Give the freshly inserted columns a sensible width (best effort - the
#    stored width is quantized by the engine's character-width conversion).
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth    = 10.166666666666666
$ws.Columns("G:I").ColumnWidth  = 11.166666666666666
$ws.Columns("M:P").ColumnWidth  = 18.451822916666668
$ws.Columns("R:V").ColumnWidth  = 11.451822916666666

# ---------------------------------------------------------------------------
# 3. Fill in the header row (row 1) for every newly-created column, left to
#    right, so the shared-string table grows in the same order as the
#    target file.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value  = "TEXT:1"
$ws.Range("G1").Value  = "CUST.SEGMENT"
$ws.Range("H1").Value  = "REL.MANAGER"
$ws.Range("I1").Value  = "CUS.TYPE.LC"
$ws.Range("M1").Value  = "NAME.2:1"
$ws.Range("N1").Value  = "TOWN.COUNTRY:1"
$ws.Range("O1").Value  = "POST.CODE:1"
$ws.Range("P1").Value  = "CUST.OFF.PHONE:1"
$ws.Range("R1").Value  = "RELATION.CODE:1"
$ws.Range("S1").Value  = "REL.CUSTOMER:1"
$ws.Range("T1").Value  = "PERSONAL.BANK"
$ws.Range("U1").Value  = "OTHER.OFFICER:1"
$ws.Range("V1").Value  = "CRC"

# ---------------------------------------------------------------------------
# 4. Append the brand-new trailing columns (AE:AS) with their headers - these
#    sit past the old last column (now AD) and have no data in rows 2/3.
# ---------------------------------------------------------------------------
$ws.Range("AE1").Value = "INTRO.ACC.NO"
$ws.Range("AF1").Value = "INTRO.MAINT.DAT"
$ws.Range("AG1").Value = "INTRO.NAME"
$ws.Range("AH1").Value = "INTRO.BNK.NAME:1"
$ws.Range("AI1").Value = "INTRO.BR.NAME:1"
$ws.Range("AJ1").Value = "ITRO.ADD:1"
$ws.Range("AK1").Value = "INTRO.CON.NO"
$ws.Range("AL1").Value = "CP.NAME:1"
$ws.Range("AM1").Value = "CP.TITLE:1"
$ws.Range("AN1").Value = "CP.ADD:1"
$ws.Range("AO1").Value = "CP.ADD2:1"
$ws.Range("AP1").Value = "CP.PH.OFF:1"
$ws.Range("AQ1").Value = "CP.FAX.NO:1"
$ws.Range("AR1").Value = "CP.CELL.NO:1"
$ws.Range("AS1").Value = "CP.PH.RES:1"

# ---------------------------------------------------------------------------
# 5. Restore the tab selection / active cell to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("AM9").Select()
